$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.810.10"
$ws.Range("E2").Value = "  +3.14%  "

$ws.Range("D3").Value = "3.471.58"
$ws.Range("E3").Value = "  +1.42%  "

$ws.Range("E4").Value = "  -0.07%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "582.57"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +2.56%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "168.44"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +7.33%  "

$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").Value = "3.472.03"
$ws.Range("E8").Value = "  +1.44%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.564"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -1.33%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "7.25"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +0.72%  "

$ws.Range("E11").Value = "  +3.57%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.430"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +0.09%  "

$ws.Range("D13").Value = "4.066.21"
$ws.Range("E13").Value = "  +0.70%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.134"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +0.31%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "27.66"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +1.83%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.0000177"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +1.86%  "

$ws.Range("D17").Value = "65.717.76"
$ws.Range("E17").Value = "  +2.79%  "

$ws.Range("D18").Value = "3.474.00"
$ws.Range("E18").Value = "  +3.07%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "6.24"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +1.40%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "13.83"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +1.09%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "384.95"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +2.39%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "7.95"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +1.64%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +0.02%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "71.76"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -0.36%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "0.522"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -0.12%  "

$ws.Range("E26").Value = "  +1.93%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "9.93"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +2.16%  "

$ws.Range("E28").Value = "  +2.17%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -0.64%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "6.27"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +3.78%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "1.45"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +3.64%  "

$ws.Range("E32").Value = "  +1.93%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "23.39"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +1.16%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "7.32"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +5.09%  "

$ws.Range("E35").Value = "  +0.03%  "

$ws.Range("E36").Value = "  -1.39%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "160.40"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -0.02%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.897"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +9.08%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "1.87"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +1.73%  "

$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "6.70"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +4.70%  "

$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.0739"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +0.49%  "

$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "26.37"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -0.35%  "

$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "2.821.65"
$ws.Range("E43").Value = "  +1.14%  "

$ws.Range("E44").Value = "  +5.63%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "43.13"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +1.15%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "4.48"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +0.83%  "

$ws.Range("B47").Value = "dogwifhat"
$ws.Range("C47").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "2.50"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +6.50%  "

$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.0308"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +0.81%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "341.97"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +4.09%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "1.08"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +4.16%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "32.43"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +6.89%  "
